# Added new screenshot keyword
# Insert a new test step "SCREENSHOT" / "checkoutScreenShotFile" right
# after the "Check checkout content" block (row 14) on the "Add to Cart"
# sheet, pushing the "Logout from application" block down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add to Cart")

# Insert a new blank row at position 14 (shifts rows 14.. down by one,
# inheriting the formatting of the row that used to be there).
$ws.Rows(14).Insert()

# Populate the new keyword row.
$ws.Range("B14").Value = "SCREENSHOT"
$ws.Range("E14").Value = "checkoutScreenShotFile"

# Column E needs to grow to fit the new, longer value.
$ws.Columns("E").ColumnWidth = 25.46875
